$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column (D) keeps its literal text formatting so Excel
# does not auto-convert strings that look like numbers/dates.

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '37.756.57'
$ws.Range("E2").Value = '  -1.10%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.027.24'
$ws.Range("E3").Value = '  -2.12%  '

# Row 4
$ws.Range("E4").Value = '  -0.05%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '227.26'
$ws.Range("E5").Value = '  -1.73%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.607'
$ws.Range("E6").Value = '  -1.78%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '59.69'
$ws.Range("E7").Value = '  +2.66%  '

# Row 8
$ws.Range("E8").Value = '  +0.00%  '

# Row 9
$ws.Range("E9").Value = '  -0.66%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0815'
$ws.Range("E10").Value = '  +0.88%  '

# Row 11
$ws.Range("E11").Value = '  -0.06%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '14.58'
$ws.Range("E12").Value = '  -0.47%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '2.326.84'
$ws.Range("E13").Value = '  -2.10%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '21.02'
$ws.Range("E14").Value = '  +1.33%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.752'
$ws.Range("E15").Value = '  -0.15%  '

# Row 16
$ws.Range("E16").Value = '  -1.06%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.038.77'
$ws.Range("E17").Value = '  -1.53%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '37.794.63'
$ws.Range("E18").Value = '  -0.83%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.00'
$ws.Range("E19").Value = '  -4.24%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '69.70'
$ws.Range("E20").Value = '  -0.54%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.0' + [string][char]0x2083 + '0822'
$ws.Range("E21").Value = '  -1.37%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '225.42'
$ws.Range("E22").Value = '  -0.04%  '

# Row 23
$ws.Range("E23").Value = '  +0.06%  '

# Row 24
$ws.Range("E24").Value = '  -2.32%  '

# Row 25
$ws.Range("E25").Value = '  -2.58%  '

# Row 26
$ws.Range("B26").Value = 'Cosmos'
$ws.Range("C26").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.20'
$ws.Range("E26").Value = '  -1.27%  '

# Row 27
$ws.Range("B27").Value = 'Monero'
$ws.Range("C27").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '164.71'
$ws.Range("E27").Value = '  -0.75%  '

# Row 28
$ws.Range("E28").Value = '  -3.51%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '18.88'
$ws.Range("E29").Value = '  -1.48%  '

# Row 30
$ws.Range("E30").Value = '  -6.28%  '

# Row 31
$ws.Range("E31").Value = '  +1.26%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.42'
$ws.Range("E32").Value = '  -3.24%  '

# Row 33
$ws.Range("E33").Value = '  +3.33%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0601'
$ws.Range("E34").Value = '  -2.66%  '

# Row 35
$ws.Range("E35").Value = '  -3.05%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.36'
$ws.Range("E36").Value = '  +6.03%  '

# Row 37
$ws.Range("E37").Value = '  -6.16%  '

# Row 38
$ws.Range("E38").Value = '  -1.88%  '

# Row 39
$ws.Range("E39").Value = '  -0.21%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.533.99'
$ws.Range("E40").Value = '  +3.55%  '

# Row 41
$ws.Range("E41").Value = '  -1.59%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '96.57'
$ws.Range("E42").Value = '  -2.16%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '16.51'
$ws.Range("E43").Value = '  -1.78%  '

# Row 44
$ws.Range("E44").Value = '  -1.51%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0918'
$ws.Range("E45").Value = '  -3.41%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.11'
$ws.Range("E46").Value = '  -2.03%  '

# Row 47
$ws.Range("E47").Value = '  -2.37%  '

# Row 48
$ws.Range("B48").Value = 'MXToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.96'
$ws.Range("E48").Value = '  -0.27%  '

# Row 49
$ws.Range("B49").Value = 'ARBITRUM'
$ws.Range("C49").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.00'
$ws.Range("E49").Value = '  -2.64%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.09'
$ws.Range("E50").Value = '  -0.58%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.217.22'
$ws.Range("E51").Value = '  -1.93%  '
